$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 / column B (the "R40" label cell) is changed to the text "1".
# Force text (not numeric) so it stays stored as a shared string, matching
# the source workbook's string-typed cell.
$ws.Range("B11").Value = "1"
